$wb = $excel.ActiveWorkbook

# --- "Forecast Comparison" sheet: MyForecast (column D) updates ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D13").Value = 1
$wsForecast.Range("D14").Value = 1
$wsForecast.Range("D17").Value = 1

# --- "Summary" sheet: text-value updates (force text, not number/date) ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "19"
$wsSummary.Range("B13").NumberFormat = "@"
$wsSummary.Range("B13").Value = "2025-04-20"
